$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed values in rows 2-7 (A4/B4 stays the same)
$ws.Range("B2").Value = 72.398
$ws.Range("B3").Value = 1648.792
$ws.Range("B5").Value = 1833.323
$ws.Range("B6").Value = 75.20399999999999
$ws.Range("B7").Value = 210.999

# Delete the 5 "Decentral_BP_*" rows (old rows 8-12); remaining rows shift up
$ws.Range("A8:B12").EntireRow.Delete() | Out-Null

# Final id column (row 8 onward, after shifting) should already hold what were
# rows 13-24; explicitly (re)assert full row 8-19 content/values to match target
$ids = @(
    "id_DK_nan_CD_Biogas",
    "id_DK_nan_IndustryE_Biogas",
    "id_DK_nan_IndustryE_Biomass",
    "id_DK_nan_CD_Coal",
    "id_DK_nan_CD_Natgas",
    "id_DK_nan_IndustryE_Natgas",
    "id_DK_nan_CD_Oil",
    "id_DK_nan_IndustryE_Oil",
    "id_DK_nan_PV",
    "id_DK_nan_ROR",
    "id_DK_nan_WL",
    "id_DK_nan_WS"
)
$vals = @(
    11.174,
    49.789,
    1.781,
    21.34,
    124.224,
    123.332,
    928.3820000000001,
    27.331,
    3910.802,
    6.894,
    4856.446,
    2655.6
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = 8 + $i
    $ws.Range("A$row").Value = $ids[$i]
    $ws.Range("B$row").Value = $vals[$i]
}
